# Trade #30 closed at 2026-02-17 15:22:29 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.84   # Current Capital
$summary.Range("B4").Value = -0.16     # Total P&L $
$summary.Range("B5").Value = -0.11     # Total P&L %
$summary.Range("B6").Value = 30        # Total Trades
$summary.Range("B8").Value = 15        # Losing Trades
$summary.Range("B9").Value = 30        # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.84      # Capital
$status.Range("D4").Value = 30         # Trades
$status.Range("E4").Value = -0.16      # P&L $
$status.Range("F4").Value = -0.16      # P&L %
$status.Range("G4").Value = 30         # Win Rate %

# --- All Trades & MarketMaking sheets, row 31 (Trade #30) ---
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($name in $tradeSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("G31").Value = 0.75              # Exit Price
    $ws.Range("H31").Value = "CLOSED"          # Status
    $ws.Range("I31").Value = -5.0633           # P&L %
    $ws.Range("J31").Value = -0.04             # P&L $
    $ws.Range("K31").Value = 99.84             # Capital After
    $ws.Range("P31").Value = "early_exit"      # Exit Reason
    $ws.Range("Q31").Value = 0.13              # Duration (min)
}
